$d = $word.ActiveDocument

# 1. After the "${keterangan_tempus}" merge field text, append a new run
#    " - ${berakhirnya_tempus}" (Arial 10pt black, same formatting as the
#    existing run it follows).
$find = $d.Content.Find
$find.Execute("`${keterangan_tempus}", $false, $false, $false, $false, $false, `
              $true, 1, $false, "", 0) | Out-Null

if ($find.Found) {
    $r = $d.Content
    $r.Start = $find.Parent.End
    $r.End = $find.Parent.End
    $r.InsertAfter(" - `${berakhirnya_tempus}")
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 0
}

# 2. Add a lastRenderedPageBreak marker before the very last run (the
#    trailing " " paragraph right before the sectPr).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertBefore([char]12)
